$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new client number (text becomes a number in A7) and the PST label in B7
$ws.Range("A7").Value = 573185344536
$ws.Range("B7").Value = "PST"

# C7 gets the same underlined formatting used elsewhere in the sheet (like C5), left empty
$ws.Range("C7").Font.Underline = 2

# Widen column A to fit the new data (matches the diff's explicit custom width)
$ws.Columns.Item(1).ColumnWidth = 24.28

# Move the active selection to C7, matching the post-edit cursor position
$ws.Range("C7").Select()
